$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '69.050.49'
$ws.Cells.Item(2, 5).Value = '  -0.13%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.799.37'
$ws.Cells.Item(3, 5).Value = '  +1.19%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '601.37'
$ws.Cells.Item(5, 5).Value = '  -0.42%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '163.18'
$ws.Cells.Item(6, 5).Value = '  -3.85%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.796.37'
$ws.Cells.Item(7, 5).Value = '  +1.15%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.07%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.536'
$ws.Cells.Item(9, 5).Value = '  +0.01%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.34%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.76%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.458'
$ws.Cells.Item(12, 5).Value = '  -1.36%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '37.23'
$ws.Cells.Item(13, 5).Value = '  -3.16%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '0.0000245'
$ws.Cells.Item(14, 5).Value = '  -1.63%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.441.64'
$ws.Cells.Item(15, 5).Value = '  +1.38%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.807.15'
$ws.Cells.Item(16, 5).Value = '  +1.39%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '69.181.36'
$ws.Cells.Item(17, 5).Value = '  +0.10%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '7.39'
$ws.Cells.Item(18, 5).Value = '  +1.21%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'TRON'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(19, 4).Value = '0.114'
$ws.Cells.Item(19, 5).Value = '  -0.36%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Chainlink'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(20, 4).Value = '17.31'
$ws.Cells.Item(20, 5).Value = '  +1.08%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '11.25'
$ws.Cells.Item(21, 5).Value = '  +3.40%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '489.40'
$ws.Cells.Item(22, 5).Value = '  -1.05%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '0.721'
$ws.Cells.Item(23, 5).Value = '  -1.09%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '0.0000151'
$ws.Cells.Item(24, 5).Value = '  -3.93%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '84.58'
$ws.Cells.Item(25, 5).Value = '  -0.98%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -3.67%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '12.21'
$ws.Cells.Item(27, 5).Value = '  -1.52%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '10.05'
$ws.Cells.Item(28, 5).Value = '  -3.71%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.09%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.59%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '8.03'
$ws.Cells.Item(31, 5).Value = '  +0.55%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -5.83%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '3.949.68'
$ws.Cells.Item(33, 5).Value = '  +1.29%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '31.86'
$ws.Cells.Item(34, 5).Value = '  -0.64%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '3.748.86'
$ws.Cells.Item(35, 5).Value = '  +1.61%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -2.20%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '0.140'
$ws.Cells.Item(37, 5).Value = '  +5.37%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.17%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '5.90'
$ws.Cells.Item(39, 5).Value = '  +0.12%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.15%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.64%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '3.03'
$ws.Cells.Item(42, 5).Value = '  -1.26%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '48.48'
$ws.Cells.Item(43, 5).Value = '  -0.69%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '1.98'
$ws.Cells.Item(44, 5).Value = '  +0.00%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '419.90'
$ws.Cells.Item(45, 5).Value = '  -4.33%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.01%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '8.38'
$ws.Cells.Item(47, 5).Value = '  -1.15%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '2.821.88'
$ws.Cells.Item(48, 5).Value = '  +1.42%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '141.79'
$ws.Cells.Item(49, 5).Value = '  +0.01%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '39.55'
$ws.Cells.Item(50, 5).Value = '  -2.83%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '1.29'
$ws.Cells.Item(51, 5).Value = '  +4.44%  '
